$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 4500
$ws.Range("J69").Value = 5000
$ws.Range("L69").Value = 15000
$ws.Range("N69").Value = -16748
$ws.Range("H72").Value = 4500
$ws.Range("J72").Value = 5000
$ws.Range("L72").Value = 45000
$ws.Range("N72").Value = -53736
$ws.Range("H105").Value = 77899.336
$ws.Range("J105").Value = 77899.336
$ws.Range("L105").Value = 77899.336
$ws.Range("N105").Value = -84887.336
$ws.Range("H132").Value = 1051.5
$ws.Range("I132").Value = 1071.0312
$ws.Range("J132").Value = 973.375
$ws.Range("K132").Value = 3213.0936
$ws.Range("L132").Value = 2920.125
$ws.Range("M132").Value = -683.0935999999997
$ws.Range("N132").Value = -7980.125
$ws.Range("H137").Value = 1348.4615
$ws.Range("I137").Value = 956
$ws.Range("J137").Value = 1976.4
$ws.Range("K137").Value = 2868
$ws.Range("L137").Value = 5929.200000000001
$ws.Range("M137").Value = -318
$ws.Range("N137").Value = -11029.2
$ws.Range("H140").Value = 77699
$ws.Range("J140").Value = 77699
$ws.Range("L140").Value = 77699
$ws.Range("N140").Value = -88059

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1164875.1
$ws.Range("I2").Value = 1455275.2
$ws.Range("J2").Value = 3274.5
$ws.Range("K2").Value = 1455275.2
$ws.Range("L2").Value = 3274.5
$ws.Range("M2").Value = -1455162.2
$ws.Range("N2").Value = -3500.5
$ws.Range("H32").Value = 2936.7092
$ws.Range("J32").Value = 9254.111000000001
$ws.Range("L32").Value = 9254.111000000001
$ws.Range("N32").Value = -9828.111000000001
$ws.Range("H35").Value = 3584
$ws.Range("I35").Value = 3584
$ws.Range("K35").Value = 3584
$ws.Range("M35").Value = -3178
$ws.Range("H45").Value = 6430537.5
$ws.Range("I45").Value = 9001791
$ws.Range("K45").Value = 9001791
$ws.Range("M45").Value = -9001414
$ws.Range("H61").Value = 3085
$ws.Range("I61").Value = 1027.2858
$ws.Range("J61").Value = 5485.6665
$ws.Range("K61").Value = 1027.2858
$ws.Range("L61").Value = 5485.6665
$ws.Range("M61").Value = -815.2858000000001
$ws.Range("N61").Value = -5909.6665
$ws.Range("H116").Value = 1164875.1
$ws.Range("I116").Value = 1455275.2
$ws.Range("J116").Value = 3274.5
$ws.Range("K116").Value = 1455275.2
$ws.Range("L116").Value = 3274.5
$ws.Range("M116").Value = -1452981.2
$ws.Range("N116").Value = -7862.5
$ws.Range("H122").Value = 3078
$ws.Range("I122").Value = 3604
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 10812
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -8362
$ws.Range("N122").Value = -9400
$ws.Range("H136").Value = 3085
$ws.Range("I136").Value = 1027.2858
$ws.Range("J136").Value = 5485.6665
$ws.Range("K136").Value = 3081.8574
$ws.Range("L136").Value = 16456.9995
$ws.Range("M136").Value = -531.8574000000003
$ws.Range("N136").Value = -21556.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1164875.1
$ws.Range("I3").Value = 1455275.2
$ws.Range("J3").Value = 3274.5
$ws.Range("K3").Value = 1455275.2
$ws.Range("L3").Value = 3274.5
$ws.Range("M3").Value = -1455161.2
$ws.Range("N3").Value = -3502.5
$ws.Range("H99").Value = 1046.1538
$ws.Range("I99").Value = 1072.8182
$ws.Range("J99").Value = 899.5
$ws.Range("K99").Value = 1072.8182
$ws.Range("L99").Value = 899.5
$ws.Range("M99").Value = 425.1818000000001
$ws.Range("N99").Value = -3895.5
$ws.Range("H107").Value = 4304.8887
$ws.Range("I107").Value = 2187.5
$ws.Range("K107").Value = 2187.5
$ws.Range("M107").Value = -267.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 497.7
$ws.Range("I107").Value = 354.91666
$ws.Range("J107").Value = 711.875
$ws.Range("K107").Value = 354.91666
$ws.Range("L107").Value = 711.875
$ws.Range("M107").Value = 1565.08334
$ws.Range("N107").Value = -4551.875
$ws.Range("H109").Value = 49900
$ws.Range("J109").Value = 49900
$ws.Range("L109").Value = 49900
$ws.Range("N109").Value = -51980
$ws.Range("H132").Value = 3553.0967
$ws.Range("I132").Value = 2856.15
$ws.Range("J132").Value = 4820.273
$ws.Range("K132").Value = 8568.450000000001
$ws.Range("L132").Value = 14460.819
$ws.Range("M132").Value = -6038.450000000001
$ws.Range("N132").Value = -19520.819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2445.5264
$ws.Range("J68").Value = 2847.7856
$ws.Range("L68").Value = 8543.356800000001
$ws.Range("N68").Value = -10165.3568
$ws.Range("H71").Value = 2445.5264
$ws.Range("J71").Value = 2847.7856
$ws.Range("L71").Value = 25630.0704
$ws.Range("N71").Value = -33742.0704
$ws.Range("H107").Value = 1811.5883
$ws.Range("J107").Value = 1813.8889
$ws.Range("L107").Value = 5441.6667
$ws.Range("N107").Value = -9281.6667
$ws.Range("H131").Value = 12821765
$ws.Range("J131").Value = 1292.2972
$ws.Range("L131").Value = 3876.8916
$ws.Range("N131").Value = -13956.8916
$ws.Range("H137").Value = 2919.276
$ws.Range("J137").Value = 4222.9443
$ws.Range("L137").Value = 12668.8329
$ws.Range("N137").Value = -22868.8329

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 904.1667
$ws.Range("I107").Value = 211
$ws.Range("K107").Value = 211
$ws.Range("M107").Value = 1709
$ws.Range("H126").Value = 1666088.6
$ws.Range("I126").Value = 2139545
$ws.Range("J126").Value = 127355.125
$ws.Range("K126").Value = 6418635
$ws.Range("L126").Value = 382065.375
$ws.Range("M126").Value = -6416165
$ws.Range("N126").Value = -387005.375
$ws.Range("H132").Value = 1604787.5
$ws.Range("I132").Value = 1674039.1
$ws.Range("K132").Value = 5022117.300000001
$ws.Range("M132").Value = -5019587.300000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 19000
$ws.Range("J63").Value = 19000
$ws.Range("L63").Value = 19000
$ws.Range("N63").Value = -20498
$ws.Range("H66").Value = 19000
$ws.Range("J66").Value = 19000
$ws.Range("L66").Value = 57000
$ws.Range("N66").Value = -64488
$ws.Range("H68").Value = 2991.5
$ws.Range("I68").Value = 2704.5715
$ws.Range("K68").Value = 2704.5715
$ws.Range("M68").Value = -1955.5715
$ws.Range("H71").Value = 2991.5
$ws.Range("I71").Value = 2704.5715
$ws.Range("K71").Value = 13522.8575
$ws.Range("M71").Value = -9778.8575
$ws.Range("H82").Value = 5035.6
$ws.Range("I82").Value = 2000
$ws.Range("J82").Value = 5794.5
$ws.Range("K82").Value = 2000
$ws.Range("L82").Value = 5794.5
$ws.Range("M82").Value = -1639
$ws.Range("N82").Value = -6516.5
$ws.Range("H85").Value = 5035.6
$ws.Range("I85").Value = 2000
$ws.Range("J85").Value = 5794.5
$ws.Range("K85").Value = 2000
$ws.Range("L85").Value = 5794.5
$ws.Range("M85").Value = -752
$ws.Range("N85").Value = -8290.5
$ws.Range("H132").Value = 3490.56
$ws.Range("J132").Value = 4706.875
$ws.Range("L132").Value = 14120.625
$ws.Range("N132").Value = -19180.625
$ws.Range("H136").Value = 4582.893
$ws.Range("I136").Value = 3282
$ws.Range("J136").Value = 6924.5
$ws.Range("K136").Value = 9846
$ws.Range("L136").Value = 20773.5
$ws.Range("M136").Value = -7296
$ws.Range("N136").Value = -25873.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 723.5
$ws.Range("I107").Value = 650
$ws.Range("J107").Value = 738.2
$ws.Range("K107").Value = 1950
$ws.Range("L107").Value = 2214.6
$ws.Range("M107").Value = -30
$ws.Range("N107").Value = -6054.6
$ws.Range("H109").Value = 42900
$ws.Range("J109").Value = 42900
$ws.Range("L109").Value = 42900
$ws.Range("N109").Value = -45674
$ws.Range("H122").Value = 172038.73
$ws.Range("I122").Value = 189092.6
$ws.Range("K122").Value = 567277.8
$ws.Range("M122").Value = -564827.8
